$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.893344666666666
$ws.Range("H2").Value = 5.680033999999999
$ws.Range("I2").Value = 0.05525983881677096
$ws.Range("J2").Value = 0.05525983881677096
$ws.Range("O2").Value = 0.02773017886769741
$ws.Range("P2").Value = 0.02773017886769741
$ws.Range("Q2").Value = 0.09834284644622221
$ws.Range("R2").Value = 0.8850856180159998
$ws.Range("S2").Value = 0.001532365214589187
$ws.Range("T2").Value = 0.001532365214589187

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.893344666666666
$ws.Range("H3").Value = 5.680033999999999
$ws.Range("I3").Value = 0.05525983881677096
$ws.Range("J3").Value = 0.05525983881677096
$ws.Range("M3").Value = 1.821156333333333
$ws.Range("N3").Value = 5.463469
$ws.Range("O3").Value = 0.9722698211323025
$ws.Range("P3").Value = 0.9722698211323026
$ws.Range("Q3").Value = 3.448076630882889
$ws.Range("R3").Value = 31.032689677946
$ws.Range("S3").Value = 0.05372747360218177
$ws.Range("T3").Value = 0.05372747360218177

# Row 4
$ws.Range("G4").Value = 4.159773333333334
$ws.Range("I4").Value = 0.1214086415227279
$ws.Range("J4").Value = 0.1214086415227279
$ws.Range("O4").Value = 0.02773017886769741
$ws.Range("P4").Value = 0.02773017886769741
$ws.Range("Q4").Value = 0.2160641732977778
$ws.Range("S4").Value = 0.0033666833455094
$ws.Range("T4").Value = 0.0033666833455094

# Row 5
$ws.Range("G5").Value = 4.159773333333334
$ws.Range("I5").Value = 0.1214086415227279
$ws.Range("J5").Value = 0.1214086415227279
$ws.Range("M5").Value = 1.821156333333333
$ws.Range("N5").Value = 5.463469
$ws.Range("O5").Value = 0.9722698211323025
$ws.Range("P5").Value = 0.9722698211323026
$ws.Range("Q5").Value = 7.575597551231112
$ws.Range("R5").Value = 68.18037796108001
$ws.Range("S5").Value = 0.1180419581772185
$ws.Range("T5").Value = 0.1180419581772185

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.240212333333333
$ws.Range("H6").Value = 3.720637
$ws.Range("I6").Value = 0.03619728348733726
$ws.Range("J6").Value = 0.03619728348733727
$ws.Range("O6").Value = 0.02773017886769741
$ws.Range("P6").Value = 0.02773017886769741
$ws.Range("Q6").Value = 0.06441828220977777
$ws.Range("R6").Value = 0.5797645398879999
$ws.Range("S6").Value = 0.001003757145628612
$ws.Range("T6").Value = 0.001003757145628613

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.240212333333333
$ws.Range("H7").Value = 3.720637
$ws.Range("I7").Value = 0.03619728348733726
$ws.Range("J7").Value = 0.03619728348733727
$ws.Range("M7").Value = 1.821156333333333
$ws.Range("N7").Value = 5.463469
$ws.Range("O7").Value = 0.9722698211323025
$ws.Range("P7").Value = 0.9722698211323026
$ws.Range("Q7").Value = 2.258620545528111
$ws.Range("R7").Value = 20.327584909753
$ws.Range("S7").Value = 0.03519352634170865
$ws.Range("T7").Value = 0.03519352634170866

# Row 8
$ws.Range("G8").Value = 26.96925
$ws.Range("H8").Value = 80.90774999999999
$ws.Range("I8").Value = 0.7871342361731639
$ws.Range("J8").Value = 0.7871342361731638
$ws.Range("O8").Value = 0.02773017886769741
$ws.Range("P8").Value = 0.02773017886769741
$ws.Range("Q8").Value = 1.400818804
$ws.Range("R8").Value = 12.607369236
$ws.Range("S8").Value = 0.02182737316197021
$ws.Range("T8").Value = 0.02182737316197021

# Row 9
$ws.Range("G9").Value = 26.96925
$ws.Range("H9").Value = 80.90774999999999
$ws.Range("I9").Value = 0.7871342361731639
$ws.Range("J9").Value = 0.7871342361731638
$ws.Range("M9").Value = 1.821156333333333
$ws.Range("N9").Value = 5.463469
$ws.Range("O9").Value = 0.9722698211323025
$ws.Range("P9").Value = 0.9722698211323026
$ws.Range("Q9").Value = 49.11522044275
$ws.Range("R9").Value = 442.03698398475
$ws.Range("S9").Value = 0.7653068630111937
$ws.Range("T9").Value = 0.7653068630111937
